$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 3.262296333333333
$ws.Range("H2").Value = 9.786889
$ws.Range("I2").Value = 0.01915820289899999
$ws.Range("J2").Value = 0.01915820289899999
$ws.Range("M2").Value = 145.7007446666667
$ws.Range("N2").Value = 437.1022340000001
$ws.Range("O2").Value = 0.2865937750105843
$ws.Range("P2").Value = 0.2865937750105843
$ws.Range("Q2").Value = 475.319005090003
$ws.Range("R2").Value = 4277.871045810027
$ws.Range("S2").Value = 0.005490621691243126
$ws.Range("T2").Value = 0.005490621691243126
$ws.Range("G3").Value = 3.262296333333333
$ws.Range("H3").Value = 9.786889
$ws.Range("I3").Value = 0.01915820289899999
$ws.Range("J3").Value = 0.01915820289899999
$ws.Range("O3").Value = 0.3320294904365841
$ws.Range("P3").Value = 0.3320294904365841
$ws.Range("Q3").Value = 550.6746510772236
$ws.Range("R3").Value = 4956.071859695013
$ws.Range("S3").Value = 0.006361088346235656
$ws.Range("T3").Value = 0.006361088346235656
$ws.Range("G4").Value = 3.262296333333333
$ws.Range("H4").Value = 9.786889
$ws.Range("I4").Value = 0.01915820289899999
$ws.Range("J4").Value = 0.01915820289899999
$ws.Range("M4").Value = 128.1261546666667
$ws.Range("N4").Value = 384.378464
$ws.Range("O4").Value = 0.2520245069956105
$ws.Range("P4").Value = 0.2520245069956105
$ws.Range("Q4").Value = 417.9854845731663
$ws.Range("R4").Value = 3761.869361158496
$ws.Range("S4").Value = 0.004828336640542348
$ws.Range("T4").Value = 0.004828336640542348
$ws.Range("G5").Value = 3.262296333333333
$ws.Range("H5").Value = 9.786889
$ws.Range("I5").Value = 0.01915820289899999
$ws.Range("J5").Value = 0.01915820289899999
$ws.Range("M5").Value = 65.761079
$ws.Range("N5").Value = 197.283237
$ws.Range("O5").Value = 0.1293522275572212
$ws.Range("P5").Value = 0.1293522275572212
$ws.Range("Q5").Value = 214.5321268977437
$ws.Range("R5").Value = 1930.789142079693
$ws.Range("S5").Value = 0.002478156220978862
$ws.Range("T5").Value = 0.002478156220978862
$ws.Range("I6").Value = 0.8527862647199704
$ws.Range("J6").Value = 0.8527862647199704
$ws.Range("M6").Value = 145.7007446666667
$ws.Range("N6").Value = 437.1022340000001
$ws.Range("O6").Value = 0.2865937750105843
$ws.Range("P6").Value = 0.2865937750105843
$ws.Range("Q6").Value = 21157.80488587864
$ws.Range("R6").Value = 190420.2439729077
$ws.Range("S6").Value = 0.2444032348832718
$ws.Range("T6").Value = 0.2444032348832718
$ws.Range("I7").Value = 0.8527862647199704
$ws.Range("J7").Value = 0.8527862647199704
$ws.Range("O7").Value = 0.3320294904365841
$ws.Range("P7").Value = 0.3320294904365841
$ws.Range("S7").Value = 0.2831501889262897
$ws.Range("T7").Value = 0.2831501889262897
$ws.Range("I8").Value = 0.8527862647199704
$ws.Range("J8").Value = 0.8527862647199704
$ws.Range("M8").Value = 128.1261546666667
$ws.Range("N8").Value = 384.378464
$ws.Range("O8").Value = 0.2520245069956105
$ws.Range("P8").Value = 0.2520245069956105
$ws.Range("Q8").Value = 18605.726329108
$ws.Range("R8").Value = 167451.536961972
$ws.Range("S8").Value = 0.2149230379386787
$ws.Range("T8").Value = 0.2149230379386787
$ws.Range("I9").Value = 0.8527862647199704
$ws.Range("J9").Value = 0.8527862647199704
$ws.Range("M9").Value = 65.761079
$ws.Range("N9").Value = 197.283237
$ws.Range("O9").Value = 0.1293522275572212
$ws.Range("P9").Value = 0.1293522275572212
$ws.Range("Q9").Value = 9549.436976111525
$ws.Range("R9").Value = 85944.93278500372
$ws.Range("S9").Value = 0.1103098029717303
$ws.Range("T9").Value = 0.1103098029717303
$ws.Range("G10").Value = 21.305189
$ws.Range("H10").Value = 63.915567
$ws.Range("I10").Value = 0.1251171236325075
$ws.Range("J10").Value = 0.1251171236325075
$ws.Range("M10").Value = 145.7007446666667
$ws.Range("N10").Value = 437.1022340000001
$ws.Range("O10").Value = 0.2865937750105843
$ws.Range("P10").Value = 0.2865937750105843
$ws.Range("Q10").Value = 3104.181902564075
$ws.Range("R10").Value = 27937.63712307668
$ws.Range("S10").Value = 0.03585778878030632
$ws.Range("T10").Value = 0.03585778878030632
$ws.Range("G11").Value = 21.305189
$ws.Range("H11").Value = 63.915567
$ws.Range("I11").Value = 0.1251171236325075
$ws.Range("J11").Value = 0.1251171236325075
$ws.Range("O11").Value = 0.3320294904365841
$ws.Range("P11").Value = 0.3320294904365841
$ws.Range("Q11").Value = 3596.309568457138
$ws.Range("R11").Value = 32366.78611611424
$ws.Range("S11").Value = 0.04154257480459258
$ws.Range("T11").Value = 0.04154257480459258
$ws.Range("G12").Value = 21.305189
$ws.Range("H12").Value = 63.915567
$ws.Range("I12").Value = 0.1251171236325075
$ws.Range("J12").Value = 0.1251171236325075
$ws.Range("M12").Value = 128.1261546666667
$ws.Range("N12").Value = 384.378464
$ws.Range("O12").Value = 0.2520245069956105
$ws.Range("P12").Value = 0.2520245069956105
$ws.Range("Q12").Value = 2729.751941016566
$ws.Range("R12").Value = 24567.76746914909
$ws.Range("S12").Value = 0.03153258140019156
$ws.Range("T12").Value = 0.03153258140019156
$ws.Range("G13").Value = 21.305189
$ws.Range("H13").Value = 63.915567
$ws.Range("I13").Value = 0.1251171236325075
$ws.Range("J13").Value = 0.1251171236325075
$ws.Range("M13").Value = 65.761079
$ws.Range("N13").Value = 197.283237
$ws.Range("O13").Value = 0.1293522275572212
$ws.Range("P13").Value = 0.1293522275572212
$ws.Range("Q13").Value = 1401.052216938931
$ws.Range("R13").Value = 12609.46995245038
$ws.Range("S13").Value = 0.0161841786474171
$ws.Range("T13").Value = 0.0161841786474171
$ws.Range("G14").Value = 0.500358
$ws.Range("H14").Value = 1.501074
$ws.Range("I14").Value = 0.002938408748521978
$ws.Range("J14").Value = 0.002938408748521978
$ws.Range("M14").Value = 145.7007446666667
$ws.Range("N14").Value = 437.1022340000001
$ws.Range("O14").Value = 0.2865937750105843
$ws.Range("P14").Value = 0.2865937750105843
$ws.Range("Q14").Value = 72.90253319992401
$ws.Range("R14").Value = 656.1227987993161
$ws.Range("S14").Value = 0.0008421296557630402
$ws.Range("T14").Value = 0.0008421296557630402
$ws.Range("G15").Value = 0.500358
$ws.Range("H15").Value = 1.501074
$ws.Range("I15").Value = 0.002938408748521978
$ws.Range("J15").Value = 0.002938408748521978
$ws.Range("O15").Value = 0.3320294904365841
$ws.Range("P15").Value = 0.3320294904365841
$ws.Range("Q15").Value = 84.460281626888
$ws.Range("R15").Value = 760.1425346419921
$ws.Range("S15").Value = 0.0009756383594661533
$ws.Range("T15").Value = 0.0009756383594661533
$ws.Range("G16").Value = 0.500358
$ws.Range("H16").Value = 1.501074
$ws.Range("I16").Value = 0.002938408748521978
$ws.Range("J16").Value = 0.002938408748521978
$ws.Range("M16").Value = 128.1261546666667
$ws.Range("N16").Value = 384.378464
$ws.Range("O16").Value = 0.2520245069956105
$ws.Range("P16").Value = 0.2520245069956105
$ws.Range("Q16").Value = 64.10894649670401
$ws.Range("R16").Value = 576.980518470336
$ws.Range("S16").Value = 0.0007405510161978403
$ws.Range("T16").Value = 0.0007405510161978403
$ws.Range("G17").Value = 0.500358
$ws.Range("H17").Value = 1.501074
$ws.Range("I17").Value = 0.002938408748521978
$ws.Range("J17").Value = 0.002938408748521978
$ws.Range("M17").Value = 65.761079
$ws.Range("N17").Value = 197.283237
$ws.Range("O17").Value = 0.1293522275572212
$ws.Range("P17").Value = 0.1293522275572212
$ws.Range("Q17").Value = 32.90408196628199
$ws.Range("R17").Value = 296.136737696538
$ws.Range("S17").Value = 0.0003800897170949445
$ws.Range("T17").Value = 0.0003800897170949445
